$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: new bold header row (mirrors the "title" style used by the
#     other header rows in the sheet, e.g. row 11 / row 19) ---
$ws.Range("B23:D23").Font.Bold = $true
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"

# --- Row 24: Micro / <10 ---
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"
$ws.Range("C24").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "General"

# --- Row 25: Small / <50 ---
$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"
$ws.Range("C25").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "General"

# --- Row 26: Medium / <250 ---
$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "<250"
$ws.Range("C26").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "General"

# --- Row 27: Large / >249 ---
$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">249"
$ws.Range("C27").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "General"
